# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for the changed rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 4
